# Add a new "Slovakia" market sheet, cloned from the existing "Portugal"
# sheet (same layout/column widths/styles/merged cells), then localize its
# two free-text cells and make it the active sheet - mirroring how the
# author produced this worksheet in the commit.

$wb = $excel.ActiveWorkbook

# Portugal is the template for the new market sheet.
$wsPortugal = $wb.Worksheets.Item("Portugal")

# Copy it to the end of the tab strip (After the last existing sheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPortugal.Copy([System.Reflection.Missing]::Value, $lastSheet)

# The copy becomes the last sheet; rename it.
$wsSlovakia = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSlovakia.Name = "Slovakia"

# Localize the market name and the accessory/model code for this sheet.
$wsSlovakia.Range("B2").Value = "Slovakia Market"
$wsSlovakia.Range("B4").Value = "NGC-2930/T3236/T3235"

# Reset the source sheet's selection to a "whole sheet" selection (as if
# the user had clicked the select-all corner before leaving it) and make
# the new Slovakia sheet the active / tab-selected one with B4 selected.
$wsPortugal.Activate()
$wsPortugal.Cells.Select()

$wsSlovakia.Activate()
$wsSlovakia.Range("B4").Select()
